$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Add a new drill as row 4 of the "drills list" table (Nr = 3):
#   "Combination Passing & Shooting Drill"
# ---------------------------------------------------------------------------

# NB: values are assigned in the same order the original author typed them
# (so new entries land in the shared-string table in the matching order).
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "Combination Passing & Shooting Drill"
$ws.Range("C4").Value = "Players line up at the starting cone. Player A (1) passes forward to Player B (2) at the first cone. Player B immediately plays the ball out wide to Player C (3) positioned on the side cone. Player C takes a positive touch, drives towards goal, and shoots. After each action, players rotate by following their pass: passer moves to the next cone, shooter returns to join the back of the starting line. The next round repeats the same pattern on the opposite side, alternating right and left."
$ws.Range("D4").Value = "- Quality of passing (accuracy, pace, and timing).`n- Movement off the ball (quick support, awareness of rotation).`n- First touch into space to set up the shot.`n- Shooting technique: composure, striking with accuracy and power.`n- Encourage players to scan and adjust body position before receiving."
$ws.Range("E4").Value = "https://youtube.com/shorts/Seib5Ps9gNE"
$ws.Hyperlinks.Add($ws.Range("E4"), "https://youtube.com/shorts/Seib5Ps9gNE")
$ws.Range("F4").Value = "Half pitch. One cone as the starting point, one cone 6–8m ahead, two side cones staggered wider towards the edge of the box, and a goal with GK. Mannequins/poles can be placed centrally to simulate defenders."
$ws.Range("H4").Value = "6–12 players + goalkeeper."
$ws.Range("I4").Value = "- Restrict to one-touch passing between cones.`n- Shooter must finish with weaker foot.`n- Add passive defender(s) to pressure shooter.`n- Introduce competition between sides (first to 5 goals)."
$ws.Range("G4").Value = "Improve combination play, forward passing, movement off the ball, and finishing under minimal pressure"

$ws.Rows.Item(4).RowHeight = 15.75

# ---------------------------------------------------------------------------
# Grow the "Drills_list" table / AutoFilter range to cover the new row
# (and one extra blank row, matching the author's manual resize to A1:I5).
# ---------------------------------------------------------------------------

$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:I5"))

# Extend the data validation rules (dropdowns) down to row 5 as well, keeping
# the original rule order (Objective/Skill list first, then free Description)
# and the original alert settings (no input/error alert boxes).
$ws.Range("G2:G5").Validation.Delete()
$ws.Range("G2:G5").Validation.Add(3, 1, 1, '"Warm-Up,Dribbling / Ball Control,Passing & Receiving,Shooting / Finishing,1v1 / Defending,Small-Sided Games,Coordination & Agility,Cool Down / Fun Game,Attacking"')
$ws.Range("G2:G5").Validation.ShowInput = $false
$ws.Range("G2:G5").Validation.ShowError = $false

$ws.Range("D2:D5").Validation.Delete()
$ws.Range("D2:D5").Validation.Add(0, 0, 0)
$ws.Range("D2:D5").Validation.InCellDropdown = $false
$ws.Range("D2:D5").Validation.ShowInput = $false
$ws.Range("D2:D5").Validation.ShowError = $false

# Keep the hidden _xlnm._FilterDatabase name in sync with the resized table.
for ($i = 1; $i -le $wb.Names.Count; $i++) {
    $n = $wb.Names.Item($i)
    if ($n.Name -like "*_FilterDatabase*") {
        $n.RefersTo = "='drills list'!`$A`$1:`$I`$5"
    }
}

# Reflect the author's final on-screen selection after entering the new row.
[void]$ws.Range("A2:I5").Select()
